# Trade #47 closed at 2026-02-17 13:28:16 - unknown UNKNOWN +0.000%
#
# Updates the Summary, Strategy Status, All Trades and MarketMaking sheets
# to reflect the newly closed trade #47.

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("Summary")
$strategyStatus = $wb.Worksheets.Item("Strategy Status")
$allTrades = $wb.Worksheets.Item("All Trades")
$marketMaking = $wb.Worksheets.Item("MarketMaking")

# --- Summary sheet updates ---
$summary.Range("B3").Value = 1197.51   # Current Capital
$summary.Range("B4").Value = -2.48     # Total P&L $
$summary.Range("B5").Value = -1.06     # Total P&L %
$summary.Range("B6").Value = 47        # Total Trades
$summary.Range("B7").Value = 19        # Winning Trades
$summary.Range("B9").Value = 40.43     # Win Rate %

# --- Strategy Status sheet updates (MarketMaking row) ---
$strategyStatus.Range("C4").Value = 97.51   # Capital
$strategyStatus.Range("D4").Value = 47      # Trades
$strategyStatus.Range("E4").Value = -2.48   # P&L $
$strategyStatus.Range("F4").Value = -2.49   # P&L %
$strategyStatus.Range("G4").Value = 40.43   # Win Rate %

# --- Append new trade (row 48) to "All Trades" and "MarketMaking" sheets ---
function Add-Trade47Row($ws) {
    $ws.Cells.Item(48, 1).Value = 47

    # Force the Date column to stay text (otherwise Excel auto-converts
    # the "2026-02-17" literal into a date serial number).
    $ws.Cells.Item(48, 2).Value = "'2026-02-17"

    $ws.Cells.Item(48, 3).Value = "13:28:10"
    $ws.Cells.Item(48, 4).Value = "MarketMaking"
    $ws.Cells.Item(48, 5).Value = "UP"
    $ws.Cells.Item(48, 6).Value = 0.92
    $ws.Cells.Item(48, 7).Value = 0.94
    $ws.Cells.Item(48, 8).Value = "CLOSED"
    $ws.Cells.Item(48, 9).Value = 2.1739
    $ws.Cells.Item(48, 10).Value = 0.02
    $ws.Cells.Item(48, 11).Value = 97.51
    $ws.Cells.Item(48, 12).Value = 0
    $ws.Cells.Item(48, 13).Value = 0
    $ws.Cells.Item(48, 14).Value = 0.6
    $ws.Cells.Item(48, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(48, 16).Value = "early_exit"
    $ws.Cells.Item(48, 17).Value = 0.13
}

Add-Trade47Row $allTrades
Add-Trade47Row $marketMaking
